$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.025.07'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '2.304.34'
$ws.Range('E3').Value = '  -0.78%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.64'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.26'
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('E7').Value = '  +2.10%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.510'
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('E10').Value = '  -0.07%  '
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '17.87'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.81'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').Value = '2.663.14'
$ws.Range('D16').Value = '2.307.90'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.781'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '42.996.80'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.66'
$ws.Range('E19').Value = '  -3.43%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.12'
$ws.Range('E21').Value = '  -1.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.29'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.14'
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.22'
$ws.Range('E28').Value = '  -0.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.93'
$ws.Range('E29').Value = '  -0.99%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.09'
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.18'
$ws.Range('E32').Value = '  -3.11%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.73'
$ws.Range('E36').Value = '  -0.78%  '
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0689'
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('E42').Value = '  +0.70%  '
$ws.Range('D43').Value = '2.000.99'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('E45').Value = '  -2.83%  '
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.57'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  -2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.71'
$ws.Range('E49').Value = '  -2.34%  '
$ws.Range('D50').Value = '2.528.84'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.75'
$ws.Range('E51').Value = '  -5.55%  '
